# Generate Report for Handback
#
# Refresh the handoff/handback timestamps on the per-language status sheets
# to reflect the latest report run:
#   - zh-cn: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
#   - de-de: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("E2").Value = "2016-03-13 05:04:46"
$wsZhCn.Range("H2").Value = "2016-03-13 05:05:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("E2").Value = "2016-03-13 05:04:49"
$wsDeDe.Range("H2").Value = "2016-03-13 05:05:09"
